# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4537
$ws1.Range("F3").Value = 2495
$ws1.Range("F10").Value = 172
$ws1.Range("F12").Value = 1697
$ws1.Range("F14").Value = 3714
$ws1.Range("F16").Value = 248

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4537
$ws4.Range("F3").Value = 2495
$ws4.Range("F12").Value = 172
$ws4.Range("F16").Value = 1697
$ws4.Range("F18").Value = 3714
$ws4.Range("F20").Value = 248
